# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 5
$ws1.Range("F5").Value = 453
$ws1.Range("F8").Value = 256
$ws1.Range("F9").Value = 14231
$ws1.Range("F12").Value = 5714
$ws1.Range("F14").Value = 64
$ws1.Range("F16").Value = 56
$ws1.Range("F17").Value = 1231
$ws1.Range("F18").Value = 5
$ws1.Range("F21").Value = 778
$ws1.Range("F22").Value = 2921
$ws1.Range("F24").Value = 10499
$ws1.Range("F25").Value = 1197
$ws1.Range("F28").Value = 3725

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 5
$ws4.Range("F6").Value = 453
$ws4.Range("F9").Value = 256
$ws4.Range("F10").Value = 14231
$ws4.Range("F13").Value = 5714
$ws4.Range("F15").Value = 64
$ws4.Range("F17").Value = 56
$ws4.Range("F18").Value = 1231
$ws4.Range("F19").Value = 5
$ws4.Range("F22").Value = 778
$ws4.Range("F23").Value = 2921
$ws4.Range("F26").Value = 10499
$ws4.Range("F27").Value = 1197
$ws4.Range("F30").Value = 3725
